$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings remain stored as text (matches source formatting)
$textCells = @('D5', 'D6', 'D10', 'D15', 'D20', 'D21', 'D22', 'D25', 'D27', 'D29', 'D30', 'D33', 'D35', 'D37', 'D38', 'D42', 'D43', 'D45', 'D47', 'D48', 'D50', 'D51')
foreach ($tc in $textCells) { $ws.Range($tc).NumberFormat = "@" }

# Apply updated cryptocurrency market data
$ws.Range('D2').Value = '58.177.16'
$ws.Range('E2').Value = '  -2.64%  '
$ws.Range('D3').Value = '3.133.78'
$ws.Range('E3').Value = '  -4.28%  '
$ws.Range('D5').Value = '524.17'
$ws.Range('E5').Value = '  -5.35%  '
$ws.Range('D6').Value = '134.17'
$ws.Range('E6').Value = '  -4.56%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.131.82'
$ws.Range('E8').Value = '  -4.41%  '
$ws.Range('E9').Value = '  -4.48%  '
$ws.Range('D10').Value = '7.24'
$ws.Range('E10').Value = '  -7.21%  '
$ws.Range('E11').Value = '  -8.13%  '
$ws.Range('E12').Value = '  -6.44%  '
$ws.Range('D13').Value = '3.670.77'
$ws.Range('E13').Value = '  -4.28%  '
$ws.Range('E14').Value = '  -0.74%  '
$ws.Range('D15').Value = '25.42'
$ws.Range('E15').Value = '  -4.42%  '
$ws.Range('D16').Value = '3.134.46'
$ws.Range('E16').Value = '  -4.13%  '
$ws.Range('D17').Value = '58.177.03'
$ws.Range('E17').Value = '  -2.83%  '
$ws.Range('E18').Value = '  -6.32%  '
$ws.Range('E19').Value = '  -5.34%  '
$ws.Range('D20').Value = '12.98'
$ws.Range('E20').Value = '  -5.21%  '
$ws.Range('D21').Value = '7.88'
$ws.Range('E21').Value = '  -6.72%  '
$ws.Range('D22').Value = '343.44'
$ws.Range('E22').Value = '  -7.70%  '
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('E24').Value = '  -4.19%  '
$ws.Range('D25').Value = '67.63'
$ws.Range('E25').Value = '  -7.01%  '
$ws.Range('D26').Value = '3.263.53'
$ws.Range('E26').Value = '  -4.25%  '
$ws.Range('D27').Value = '0.171'
$ws.Range('E27').Value = '  +0.44%  '
$ws.Range('D28').Value = '0.0₃0950'
$ws.Range('E28').Value = '  -5.77%  '
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.89%  '
$ws.Range('D30').Value = '6.83'
$ws.Range('E30').Value = '  -2.83%  '
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('E32').Value = '  -7.48%  '
$ws.Range('D33').Value = '6.89'
$ws.Range('E33').Value = '  -7.15%  '
$ws.Range('E34').Value = '  +0.90%  '
$ws.Range('D35').Value = '21.34'
$ws.Range('E35').Value = '  -4.76%  '
$ws.Range('E36').Value = '  -4.03%  '
$ws.Range('D37').Value = '157.47'
$ws.Range('E37').Value = '  -5.32%  '
$ws.Range('D38').Value = '6.23'
$ws.Range('E38').Value = '  -5.29%  '
$ws.Range('E39').Value = '  -9.97%  '
$ws.Range('E40').Value = '  -5.03%  '
$ws.Range('D41').Value = '3.165.12'
$ws.Range('E41').Value = '  -4.22%  '
$ws.Range('D42').Value = '40.45'
$ws.Range('E42').Value = '  -2.29%  '
$ws.Range('D43').Value = '23.80'
$ws.Range('E43').Value = '  -8.02%  '
$ws.Range('E44').Value = '  -1.55%  '
$ws.Range('D45').Value = '0.690'
$ws.Range('E45').Value = '  -7.11%  '
$ws.Range('E46').Value = '  -4.23%  '
$ws.Range('D47').Value = '1.00'
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = '1.44'
$ws.Range('E48').Value = '  -7.26%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.279.70'
$ws.Range('E49').Value = '  -1.20%  '
$ws.Range('D50').Value = '6.14'
$ws.Range('E50').Value = '  -2.55%  '
$ws.Range('D51').Value = '20.67'
$ws.Range('E51').Value = '  -1.72%  '
